# Fruta / hortaliza, semanal
# Weekly data refresh: a new (most-recent) week's entry is inserted, the
# previously-recorded daily rows shift their dates back one slot, and the
# final "Especial"/"Primera" pair (previously rows 8-9) is re-split across
# rows 9-10 so a brand new row 10 is created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new values (new latest-week "Primera" entry) ---
$ws.Range("D4").Value = 44487
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 23000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 23500
$ws.Range("S4").Value = 2350

# --- Row 5: date shifts back one slot, rest unchanged ---
$ws.Range("D5").Value = 44446

# --- Row 6: date shifts back one slot, rest unchanged ---
$ws.Range("D6").Value = 44447

# --- Row 7: date shifts back one slot, rest unchanged ---
$ws.Range("D7").Value = 44448

# --- Row 8: date shifts back one slot; quality/prices become the old row7 values ---
$ws.Range("D8").Value = 44452
$ws.Range("L8").Value = "Primera"
$ws.Range("N8").Value = 21000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21500
$ws.Range("S8").Value = 2150

# --- Row 9: same date, becomes the "Especial" entry (old row8 values) ---
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 31000
$ws.Range("O9").Value = 32000
$ws.Range("P9").Value = 31500
$ws.Range("S9").Value = 3150

# --- Row 10: brand new row (the old "Primera" entry that used to live in row 9) ---
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D10").Value = 44461
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = "Otros"
$ws.Range("I10").Value = 100107002
$ws.Range("J10").Value = "Chirimoya"
$ws.Range("K10").Value = "Cultivar IV Región"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 30000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 30000
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 3000
$ws.Range("T10").Value = 10
